$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as plain text, even if it looks like a
# number (e.g. "303.96"), and without altering the cell's style (no NumberFormat
# / quote-prefix side effects). We do this by writing a text formula that
# evaluates to the desired literal string, then converting the cell to a
# static value via copy / paste-special-values (this mirrors what Excel does
# when you paste "Values Only", and avoids Excel's automatic text->number
# coercion that a direct .Value assignment would trigger).
function Set-TextValue($range, $text) {
    $range.Formula = '=""&"' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "41.804.75"
Set-TextValue $ws.Range("E2") "  -0.42%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.255.42"
Set-TextValue $ws.Range("E3") "  -0.75%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  +0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "303.96"
Set-TextValue $ws.Range("E5") "  -0.67%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "94.14"
Set-TextValue $ws.Range("E6") "  +1.05%  "

# Row 7 - XRP
Set-TextValue $ws.Range("E7") "  -1.23%  "

# Row 8 - USDC
Set-TextValue $ws.Range("E8") "  +0.02%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.486"
Set-TextValue $ws.Range("E9") "  -0.39%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "34.64"
Set-TextValue $ws.Range("E10") "  +5.38%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("E11") "  -1.75%  "

# Row 12 - TRON
Set-TextValue $ws.Range("E12") "  -0.37%  "

# Row 13 - Polkadot
Set-TextValue $ws.Range("D13") "6.60"

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("E14") "  -0.83%  "

# Row 15 - Chainlink
Set-TextValue $ws.Range("D15") "14.28"
Set-TextValue $ws.Range("E15") "  -0.62%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.269.25"

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.786"
Set-TextValue $ws.Range("E17") "  +0.23%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "41.707.47"
Set-TextValue $ws.Range("E18") "  -0.41%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D19") "12.28"
Set-TextValue $ws.Range("E19") "  -3.31%  "

# Row 20 - ShibaInu
Set-TextValue $ws.Range("E20") "  -2.11%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "5.93"
Set-TextValue $ws.Range("E21") "  -0.98%  "

# Row 22 - Litecoin
Set-TextValue $ws.Range("D22") "67.92"
Set-TextValue $ws.Range("E22") "  -0.40%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "236.62"
Set-TextValue $ws.Range("E23") "  -3.06%  "

# Row 24 - PancakeSwap
Set-TextValue $ws.Range("E24") "  -1.43%  "

# Row 25 - Dai
Set-TextValue $ws.Range("E25") "  -0.08%  "

# Row 26 - ImmutableX
Set-TextValue $ws.Range("E26") "  -1.58%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "23.56"
Set-TextValue $ws.Range("E27") "  -1.97%  "

# Row 28 - InjectiveProtocol
Set-TextValue $ws.Range("D28") "36.12"
Set-TextValue $ws.Range("E28") "  +2.89%  "

# Row 29 - now Cosmos (was Toncoin)
Set-TextValue $ws.Range("B29") "Cosmos"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D29") "9.41"
Set-TextValue $ws.Range("E29") "  -3.10%  "

# Row 30 - now Toncoin (was Cosmos)
Set-TextValue $ws.Range("B30") "Toncoin"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "2.10"
Set-TextValue $ws.Range("E30") "  +0.76%  "

# Row 31 - Monero
Set-TextValue $ws.Range("D31") "159.73"
Set-TextValue $ws.Range("E31") "  +0.38%  "

# Row 32 - now FirstDigitalUSD (was Filecoin)
Set-TextValue $ws.Range("B32") "FirstDigitalUSD"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D32") "1.00"
Set-TextValue $ws.Range("E32") "  +0.05%  "

# Row 33 - now Filecoin (was FirstDigitalUSD)
Set-TextValue $ws.Range("B33") "Filecoin"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "5.19"
Set-TextValue $ws.Range("E33") "  -3.47%  "

# Row 34 - LidoDAOToken
Set-TextValue $ws.Range("E34") "  +3.54%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0731"
Set-TextValue $ws.Range("E35") "  -1.91%  "

# Row 36 - Celestia
Set-TextValue $ws.Range("D36") "16.89"
Set-TextValue $ws.Range("E36") "  -2.24%  "

# Row 37 - WEMIXToken
Set-TextValue $ws.Range("E37") "  +0.55%  "

# Row 38 - Kaspa
Set-TextValue $ws.Range("E38") "  -1.13%  "

# Row 39 - ARBITRUM
Set-TextValue $ws.Range("E39") "  +0.84%  "

# Row 40 - Stellar
Set-TextValue $ws.Range("E40") "  -2.67%  "

# Row 41 - RenderToken
Set-TextValue $ws.Range("D41") "3.97"
Set-TextValue $ws.Range("E41") "  +0.33%  "

# Row 42 - ApeXProtocol
Set-TextValue $ws.Range("D42") "2.37"
Set-TextValue $ws.Range("E42") "  +5.40%  "

# Row 43 - Maker
Set-TextValue $ws.Range("D43") "1.971.92"
Set-TextValue $ws.Range("E43") "  -2.09%  "

# Row 44 - VeChain
Set-TextValue $ws.Range("D44") "0.0281"
Set-TextValue $ws.Range("E44") "  -0.65%  "

# Row 45 - EnergySwap
Set-TextValue $ws.Range("D45") "18.74"
Set-TextValue $ws.Range("E45") "  -5.15%  "

# Row 46 - NEARProtocol
Set-TextValue $ws.Range("E46") "  -0.96%  "

# Row 47 - FraxShare
Set-TextValue $ws.Range("D47") "9.81"
Set-TextValue $ws.Range("E47") "  -4.79%  "

# Row 48 - MultiversX
Set-TextValue $ws.Range("D48") "52.98"
Set-TextValue $ws.Range("E48") "  -0.89%  "

# Row 49 - BitcoinSV
Set-TextValue $ws.Range("D49") "72.56"
Set-TextValue $ws.Range("E49") "  -0.07%  "

# Row 50 - Stacks
Set-TextValue $ws.Range("E50") "  -1.21%  "

# Row 51 - Aave
Set-TextValue $ws.Range("D51") "90.48"
Set-TextValue $ws.Range("E51") "  -1.47%  "

$excel.CutCopyMode = $false
